# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stock holdings) sheet gets three new trailing columns:
#   H = date              (the filing date of this property declaration)
#   I = legislator_name   (the legislator this filing belongs to)
#   J = legislator_id     (numeric id of that legislator)
#
# Values for every existing data row (2-8):
#   date            = "2013-05-01"
#   legislator_name = "顏寬恒"
#   legislator_id   = 1803

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 8   # rows 2..8 hold the stock entries (row 1 is the header)

# ---- header row (row 1) -------------------------------------------------
$ws.Cells.Item(1, 8).Value  = "date"
$ws.Cells.Item(1, 9).Value  = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# Match the look of the existing bold/centered/bordered header cells (style
# used by B1:G1) so the new headers are visually consistent with them.
$headerRange = $ws.Range("H1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous

# ---- data rows (2..8) ----------------------------------------------------
# Force the date column to text so the ISO-looking string "2013-05-01" is
# stored verbatim instead of being auto-converted into a date serial number.
$ws.Range("H2:H" + $lastRow).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value  = "2013-05-01"
    $ws.Cells.Item($r, 9).Value  = "顏寬恒"
    $ws.Cells.Item($r, 10).Value = 1803
}
